$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text (General-looking) storage for columns B-E so that Excel
# does not auto-convert numeric- or percent-looking strings into
# numbers, matching the original workbook which stores these as text.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '319.06'
$ws.Range("E2").Value = '4.47%'
$ws.Range("E3").Value = '0.15%'
$ws.Range("D4").Value = '5.129'
$ws.Range("E4").Value = '0.36%'
$ws.Range("D5").Value = '0.08219'
$ws.Range("E5").Value = '4.47%'
$ws.Range("D6").Value = '2.152'
$ws.Range("E6").Value = '-0.98%'
$ws.Range("D7").Value = '8.007'
$ws.Range("E7").Value = '1.04%'
$ws.Range("B8").Value = 'GateToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D8").Value = '4.137'
$ws.Range("E8").Value = '0.90%'
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '0.9255'
$ws.Range("E9").Value = '0.76%'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").Value = '0.1007'
$ws.Range("E10").Value = '3.96%'
$ws.Range("B11").Value = 'WazirX'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D11").Value = '0.1895'
$ws.Range("E11").Value = '1.57%'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.09196'
$ws.Range("E12").Value = '5.78%'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03597'
$ws.Range("E13").Value = '3.22%'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09919'
$ws.Range("E14").Value = '-0.10%'
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001433'
$ws.Range("E15").Value = '-0.61%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").Value = '0.005662'
$ws.Range("E16").Value = '-0.22%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").Value = '3.462'
$ws.Range("E17").Value = '0.01%'
$ws.Range("D18").Value = '2.801'
$ws.Range("E18").Value = '17.07%'
$ws.Range("E19").Value = '-1.45%'
$ws.Range("D20").Value = '0.1331'
$ws.Range("E20").Value = '2.38%'
$ws.Range("D21").Value = '5.063'
$ws.Range("E21").Value = '4.31%'
$ws.Range("E22").Value = '-0.54%'
$ws.Range("D23").Value = '0.04589'
$ws.Range("E23").Value = '0.70%'
$ws.Range("D24").Value = '0.001246'
$ws.Range("E24").Value = '1.05%'
$ws.Range("D25").Value = '0.004738'
$ws.Range("E25").Value = '-6.90%'
$ws.Range("D26").Value = '0.0001300'
$ws.Range("E26").Value = '-7.19%'
$ws.Range("D27").Value = '0.0004499'
$ws.Range("E27").Value = '-5.30%'
$ws.Range("E39").Value = '8.89%'
$ws.Range("E40").Value = '4.43%'
$ws.Range("D41").Value = '0.007748'
$ws.Range("E41").Value = '0.38%'
$ws.Range("E42").Value = '0.13%'
$ws.Range("D43").Value = '0.007524'
$ws.Range("E43").Value = '-3.92%'
$ws.Range("D44").Value = '0.002095'
$ws.Range("E44").Value = '-6.05%'
$ws.Range("D45").Value = '0.01199'
$ws.Range("E45").Value = '5.82%'
$ws.Range("D46").Value = '0.00006455'
$ws.Range("E46").Value = '0.32%'
$ws.Range("E47").Value = '-0.04%'
$ws.Range("E48").Value = '17.98%'
$ws.Range("D49").Value = '0.001900'
$ws.Range("E49").Value = '-5.03%'
$ws.Range("D50").Value = '0.00002100'
$ws.Range("E50").Value = '-0.04%'
$ws.Range("D51").Value = '0.0002000'
$ws.Range("E51").Value = '-0.04%'

# Restore original General number format now that text values are set.
$ws.Range("B2:E51").NumberFormat = "General"

